$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 19.09702460569844
$ws.Range("R2").Value = 171.873221451286
$ws.Range("S2").Value = 0.2299953477621856
$ws.Range("T2").Value = 0.2299953477621856

# Row 3 updates
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("S3").Value = 0.6794731949692173
$ws.Range("T3").Value = 0.6794731949692174

# Row 4 updates
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 7.517027991521777
$ws.Range("R4").Value = 67.653251923696
$ws.Range("S4").Value = 0.09053145726859702
$ws.Range("T4").Value = 0.09053145726859703
